# Updates the cryptos list (Price / Volume(1h) columns) with refreshed
# quotes, mirroring the GitHub Actions data-refresh commit.
#
# Column D (Price) values that look like plain decimals (e.g. "575.46")
# would otherwise be auto-converted to numbers by Excel when assigned via
# .Value, so we force those specific cells to Text format first. Values
# that use the site's "thousands dot" style (e.g. "61.539.79") are never
# parsed as numbers and are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.539.79"
$ws.Cells.Item(2, 5).Value = "  +0.63%  "

$ws.Cells.Item(3, 4).Value = "3.390.49"
$ws.Cells.Item(3, 5).Value = "  -0.36%  "

$ws.Cells.Item(4, 5).Value = "  +0.01%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "575.46"
$ws.Cells.Item(5, 5).Value = "  +0.35%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "141.12"
$ws.Cells.Item(6, 5).Value = "  -1.00%  "

$ws.Cells.Item(7, 5).Value = "  +0.03%  "

$ws.Cells.Item(8, 5).Value = "  -0.85%  "

$ws.Cells.Item(9, 5).Value = "  +0.76%  "

$ws.Cells.Item(10, 5).Value = "  -1.46%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.385"
$ws.Cells.Item(11, 5).Value = "  -2.93%  "

$ws.Cells.Item(12, 4).Value = "3.968.61"
$ws.Cells.Item(12, 5).Value = "  -0.36%  "

$ws.Cells.Item(13, 5).Value = "  +0.24%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "28.42"
$ws.Cells.Item(14, 5).Value = "  +1.18%  "

$ws.Cells.Item(15, 4).Value = "3.389.86"
$ws.Cells.Item(15, 5).Value = "  -0.28%  "

$ws.Cells.Item(16, 5).Value = "  -0.90%  "

$ws.Cells.Item(17, 4).Value = "61.586.44"
$ws.Cells.Item(17, 5).Value = "  +0.70%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "6.12"
$ws.Cells.Item(18, 5).Value = "  -0.35%  "

$ws.Cells.Item(19, 5).Value = "  -2.06%  "

$ws.Cells.Item(20, 5).Value = "  +0.12%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "391.06"
$ws.Cells.Item(21, 5).Value = "  +1.89%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "74.98"
$ws.Cells.Item(22, 5).Value = "  +0.59%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.550"
$ws.Cells.Item(23, 5).Value = "  -1.61%  "

$ws.Cells.Item(24, 5).Value = "  +0.11%  "

$ws.Cells.Item(25, 5).Value = "  -4.82%  "

$ws.Cells.Item(26, 5).Value = "  +7.42%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.999"
$ws.Cells.Item(27, 5).Value = "  -0.07%  "

$ws.Cells.Item(28, 5).Value = "  -1.59%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.02"
$ws.Cells.Item(29, 5).Value = "  -0.22%  "

$ws.Cells.Item(30, 5).Value = "  -1.21%  "

$ws.Cells.Item(31, 5).Value = "  +0.04%  "

$ws.Cells.Item(32, 5).Value = "  -1.65%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "23.26"
$ws.Cells.Item(33, 5).Value = "  -1.18%  "

$ws.Cells.Item(34, 5).Value = "  -1.85%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "168.62"
$ws.Cells.Item(35, 5).Value = "  +0.51%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.01"
$ws.Cells.Item(36, 5).Value = "  -0.08%  "

$ws.Cells.Item(37, 4).Value = "3.424.55"
$ws.Cells.Item(37, 5).Value = "  -0.25%  "

$ws.Cells.Item(38, 5).Value = "  -1.45%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0764"
$ws.Cells.Item(39, 5).Value = "  -1.37%  "

$ws.Cells.Item(40, 5).Value = "  -5.24%  "

$ws.Cells.Item(41, 5).Value = "  -0.38%  "

$ws.Cells.Item(42, 5).Value = "  -0.69%  "

$ws.Cells.Item(43, 5).Value = "  -1.68%  "

$ws.Cells.Item(44, 5).Value = "  +2.05%  "

$ws.Cells.Item(45, 4).Value = "2.478.07"
$ws.Cells.Item(45, 5).Value = "  -0.43%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "22.88"
$ws.Cells.Item(46, 5).Value = "  -1.04%  "

$ws.Cells.Item(47, 5).Value = "  -2.67%  "

$ws.Cells.Item(48, 5).Value = "  +0.10%  "

$ws.Cells.Item(49, 5).Value = "  -1.35%  "

$ws.Cells.Item(50, 5).Value = "  -3.96%  "

$ws.Cells.Item(51, 5).Value = "  -2.39%  "
